# Update the multiplication problems/answers in the table to new values,
# per the commit "Update master to output generated at c8c62b6".

$d = $word.ActiveDocument

$replacements = @(
    @{old = "52×44=2288"; new = "25×25=625"},
    @{old = "94×96=9024"; new = "15×99=1485"},
    @{old = "87×39=3393"; new = "57×70=3990"},
    @{old = "88×16=1408"; new = "52×76=3952"},
    @{old = "75×13=975";  new = "47×77=3619"},
    @{old = "74×78=5772"; new = "91×47=4277"},
    @{old = "54×73=3942"; new = "78×80=6240"},
    @{old = "57×78=4446"; new = "72×63=4536"},
    @{old = "15×50=750";  new = "35×66=2310"},
    @{old = "26×97=2522"; new = "66×46=3036"},
    @{old = "16×57=912";  new = "96×94=9024"},
    @{old = "94×93=8742"; new = "79×29=2291"},
    @{old = "59×17=1003"; new = "54×36=1944"},
    @{old = "21×12=252";  new = "97×80=7760"},
    @{old = "25×34=850";  new = "99×43=4257"},
    @{old = "28×60=1680"; new = "11×49=539"},
    @{old = "93×81=7533"; new = "94×17=1598"},
    @{old = "54×20=1080"; new = "27×72=1944"},
    @{old = "49×19=931";  new = "72×61=4392"},
    @{old = "19×83=1577"; new = "84×88=7392"},
    @{old = "25×85=2125"; new = "33×80=2640"},
    @{old = "19×17=323";  new = "87×31=2697"},
    @{old = "44×48=2112"; new = "72×84=6048"},
    @{old = "61×79=4819"; new = "77×56=4312"},
    @{old = "82×97=7954"; new = "37×38=1406"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
